$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 13 with the match data (Charleroi 1-0 Anderlecht, 24/10/2025)
$ws.Range("A13").Value = "24/10/2025"
$ws.Range("B13").Value = "Charleroi"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "Anderlecht"
$ws.Range("F13").Value = "W"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1.13
$ws.Range("L13").Value = 1.39
$ws.Range("M13").Value = 8
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 3
$ws.Range("P13").Value = 2
